$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-computed average S-map coefficient of abundance for affected species rows
# Row 5: Merlangius merlangus
$ws.Range("C5").Value = -0.518
$ws.Range("D5").Value = 0.2966

# Row 7: Pollachius virens
$ws.Range("C7").Value = 0.3419
$ws.Range("D7").Value = 0.1978

# Update the active selection to match the recorded view state
$ws.Range("G2:I12").Select()
